$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 22, pushing the "election / presidential election"
# rule (and everything after it) down by two rows.
$ws.Rows(22).Insert()
$ws.Rows(22).Insert()

# New rule: significance / even greater significance
$ws.Cells.Item(22, 1).Value = "significance"
$ws.Cells.Item(22, 2).Value = "even greater significance"

# New rule: prelude / prelude
$ws.Cells.Item(23, 1).Value = "prelude"
$ws.Cells.Item(23, 2).Value = "prelude"

# Update the synonym list for "elections" (row 8)
$ws.Range("B8").Value = "upcoming parliamentary elections, these elections, recent elections, many elections, elections, country 's elections"

# Update the synonym list for "drive" (row 10)
$ws.Range("B10").Value = "petition drive, drive"

# Update the synonym list for "type" (originally row 29, now row 31 after the insert)
$ws.Range("B31").Value = "this type, other common type"

# New rule appended at the end: finger / finger
$ws.Cells.Item(38, 1).Value = "finger"
$ws.Cells.Item(38, 2).Value = "finger"
